# Insert 17 empty paragraphs (style "BodyText", language sr-Latn-CS) right
# after the Body Text paragraph ending in "...Predavaci za kurseve su
# nadeni i kontaktirani." and right before the "Cena realizacije projekta"
# Heading 1 paragraph.

$d = $word.ActiveDocument

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "kontaktirani" -and $p.Style.NameLocal -eq "Body Text") {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    throw "Anchor paragraph (ending in 'kontaktirani.') not found"
}

$fragment = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='BodyText'/><w:rPr><w:lang w:val='sr-Latn-CS'/></w:rPr></w:pPr></w:p>"

$pos = $anchor.Range.End

for ($i = 0; $i -lt 17; $i++) {
    $ins = $d.Range($pos, $pos)
    [void]$ins.InsertXML($fragment)
    $pos = $ins.End
}

Write-Output "Inserted 17 empty BodyText paragraphs after anchor."
